# "Ready for final reactive testing"
# Update the Diversification sheet's proposed-trade inputs (J8:J12), turn
# J12 from a formula into a plain entered value, shrink row 11's height,
# zoom all three sheets out to 80%, and leave the selection on J10 of the
# Diversification sheet (matching the author's final reactive-testing view).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Diversification")

# Proposed-trades column (J8:J12)
$ws.Range("J8").Value = 72.56
$ws.Range("J9").Value = -10
$ws.Range("J10").Value = 7
$ws.Range("J11").Value = 10
# J12 previously held "=B12+D12-E12"; replace it with a plain typed value.
$ws.Range("J12").Value = 13

# Row 11 shrinks slightly (no longer auto/default height).
$ws.Rows.Item(11).RowHeight = 19.7

# Zoom every sheet from 100% to 80%, then restore Diversification as the
# active sheet/tab with J10 selected.
$ws2 = $wb.Worksheets.Item("Sheet2")
[void]$ws2.Activate()
$excel.ActiveWindow.Zoom = 80

$ws3 = $wb.Worksheets.Item("Sheet3")
[void]$ws3.Activate()
$excel.ActiveWindow.Zoom = 80

[void]$ws.Activate()
$excel.ActiveWindow.Zoom = 80
[void]$ws.Range("J10").Select()
